$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the data range so numeric-looking strings
# (IDs, amounts, codes) are preserved exactly as text, matching the source data.
$ws.Range("A1:L15").NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = "''U56733636"
$ws.Range("B1").Value = "''9652408096438"
$ws.Range("C1").Value = "''240516090004339"
$ws.Range("D1").Value = "LIAO"
$ws.Range("E1").Value = "YONGXIANG"
$ws.Range("F1").Value = "12/11/2023"
$ws.Range("G1").Value = "900.00"
$ws.Range("H1").Value = "240.00"
$ws.Range("I1").Value = "416.70"
$ws.Range("J1").Value = "0.00"
$ws.Range("K1").Value = "0.0"
$ws.Range("L1").Value = "921279283"
$ws.Range("N1").Value = "NO"
$ws.Range("O1").Value = "NO"
$ws.Range("P1").Value = "NO"
$ws.Range("Q1").Value = "NO"

# Row 2
$ws.Range("A2").Value = "''U56733636"
$ws.Range("B2").Value = "''9652408096440"
$ws.Range("C2").Value = "''240516090004339"
$ws.Range("D2").Value = "LIAO"
$ws.Range("E2").Value = "YONGXIANG"
$ws.Range("F2").Value = "11/13/2023"
$ws.Range("G2").Value = "900.00"
$ws.Range("H2").Value = "319.81"
$ws.Range("I2").Value = "416.70"
$ws.Range("J2").Value = "0.00"
$ws.Range("K2").Value = "0.0"
$ws.Range("L2").Value = "921279283"
$ws.Range("N2").Value = "NO"
$ws.Range("O2").Value = "NO"
$ws.Range("P2").Value = "NO"
$ws.Range("Q2").Value = "NO"

# Row 3
$ws.Range("A3").Value = "''U56733636"
$ws.Range("B3").Value = "''9652408096441"
$ws.Range("C3").Value = "''240516090004339"
$ws.Range("D3").Value = "LIAO"
$ws.Range("E3").Value = "YONGXIANG"
$ws.Range("F3").Value = "11/29/2023"
$ws.Range("G3").Value = "900.00"
$ws.Range("H3").Value = "319.81"
$ws.Range("I3").Value = "416.70"
$ws.Range("J3").Value = "0.00"
$ws.Range("K3").Value = "0.0"
$ws.Range("L3").Value = "921279283"
$ws.Range("N3").Value = "NO"
$ws.Range("O3").Value = "NO"
$ws.Range("P3").Value = "NO"
$ws.Range("Q3").Value = "NO"

# Row 4
$ws.Range("A4").Value = "''U56733636"
$ws.Range("B4").Value = "''9652408096443"
$ws.Range("C4").Value = "''240516090004339"
$ws.Range("D4").Value = "LIAO"
$ws.Range("E4").Value = "YONGXIANG"
$ws.Range("F4").Value = "12/07/2023"
$ws.Range("G4").Value = "900.00"
$ws.Range("H4").Value = "319.81"
$ws.Range("I4").Value = "416.70"
$ws.Range("J4").Value = "0.00"
$ws.Range("K4").Value = "0.0"
$ws.Range("L4").Value = "921279283"
$ws.Range("N4").Value = "NO"
$ws.Range("O4").Value = "NO"
$ws.Range("P4").Value = "NO"
$ws.Range("Q4").Value = "NO"

# Row 5
$ws.Range("A5").Value = "''U56733636"
$ws.Range("B5").Value = "''9652408096663"
$ws.Range("C5").Value = "''240516090004339"
$ws.Range("D5").Value = "LIAO"
$ws.Range("E5").Value = "YONGXIANG"
$ws.Range("F5").Value = "11/20/2023"
$ws.Range("G5").Value = "900.00"
$ws.Range("H5").Value = "319.81"
$ws.Range("I5").Value = "416.70"
$ws.Range("J5").Value = "0.00"
$ws.Range("K5").Value = "0.0"
$ws.Range("L5").Value = "921279283"
$ws.Range("N5").Value = "NO"
$ws.Range("O5").Value = "NO"
$ws.Range("P5").Value = "NO"
$ws.Range("Q5").Value = "NO"

# Row 6
$ws.Range("A6").Value = "''U91200304"
$ws.Range("B6").Value = "''4682412800131"
$ws.Range("C6").Value = "''240516090004339"
$ws.Range("D6").Value = "SALTOS"
$ws.Range("E6").Value = "ELENI"
$ws.Range("F6").Value = "02/14/2024"
$ws.Range("G6").Value = "800.00"
$ws.Range("H6").Value = "65.03"
$ws.Range("I6").Value = "92.90"
$ws.Range("J6").Value = "0.00"
$ws.Range("K6").Value = "0.0"
$ws.Range("L6").Value = "921279283"
$ws.Range("N6").Value = "NO"
$ws.Range("O6").Value = "NO"
$ws.Range("P6").Value = "NO"
$ws.Range("Q6").Value = "NO"

# Row 7
$ws.Range("A7").Value = "''U91200304"
$ws.Range("B7").Value = "''4682412800132"
$ws.Range("C7").Value = "''240516090004339"
$ws.Range("D7").Value = "SALTOS"
$ws.Range("E7").Value = "ELENI"
$ws.Range("F7").Value = "02/28/2024"
$ws.Range("G7").Value = "800.00"
$ws.Range("H7").Value = "65.03"
$ws.Range("I7").Value = "92.90"
$ws.Range("J7").Value = "0.00"
$ws.Range("K7").Value = "0.0"
$ws.Range("L7").Value = "921279283"
$ws.Range("N7").Value = "NO"
$ws.Range("O7").Value = "NO"
$ws.Range("P7").Value = "NO"
$ws.Range("Q7").Value = "NO"

# Row 8
$ws.Range("A8").Value = "''U91200304"
$ws.Range("B8").Value = "''4682412800133"
$ws.Range("C8").Value = "''240516090004339"
$ws.Range("D8").Value = "SALTOS"
$ws.Range("E8").Value = "ELENI"
$ws.Range("F8").Value = "03/06/2024"
$ws.Range("G8").Value = "800.00"
$ws.Range("H8").Value = "65.03"
$ws.Range("I8").Value = "92.90"
$ws.Range("J8").Value = "0.00"
$ws.Range("K8").Value = "0.0"
$ws.Range("L8").Value = "921279283"
$ws.Range("N8").Value = "NO"
$ws.Range("O8").Value = "NO"
$ws.Range("P8").Value = "NO"
$ws.Range("Q8").Value = "NO"

# Row 9
$ws.Range("A9").Value = "''U91200304"
$ws.Range("B9").Value = "''4682412800134"
$ws.Range("C9").Value = "''240516090004339"
$ws.Range("D9").Value = "SALTOS"
$ws.Range("E9").Value = "ELENI"
$ws.Range("F9").Value = "03/27/2024"
$ws.Range("G9").Value = "800.00"
$ws.Range("H9").Value = "66.12"
$ws.Range("I9").Value = "94.46"
$ws.Range("J9").Value = "0.00"
$ws.Range("K9").Value = "0.0"
$ws.Range("L9").Value = "921279283"
$ws.Range("N9").Value = "NO"
$ws.Range("O9").Value = "NO"
$ws.Range("P9").Value = "NO"
$ws.Range("Q9").Value = "NO"

# Row 10
$ws.Range("A10").Value = "''U91200304"
$ws.Range("B10").Value = "''4682412800135"
$ws.Range("C10").Value = "''240516090004339"
$ws.Range("D10").Value = "SALTOS"
$ws.Range("E10").Value = "ELENI"
$ws.Range("F10").Value = "03/13/2024"
$ws.Range("G10").Value = "800.00"
$ws.Range("H10").Value = "66.12"
$ws.Range("I10").Value = "94.46"
$ws.Range("J10").Value = "0.00"
$ws.Range("K10").Value = "0.0"
$ws.Range("L10").Value = "921279283"
$ws.Range("N10").Value = "NO"
$ws.Range("O10").Value = "NO"
$ws.Range("P10").Value = "NO"
$ws.Range("Q10").Value = "NO"

# Row 11
$ws.Range("A11").Value = "''U91200304"
$ws.Range("B11").Value = "''4682412800136"
$ws.Range("C11").Value = "''240516090004339"
$ws.Range("D11").Value = "SALTOS"
$ws.Range("E11").Value = "ELENI"
$ws.Range("F11").Value = "04/03/2024"
$ws.Range("G11").Value = "800.00"
$ws.Range("H11").Value = "66.12"
$ws.Range("I11").Value = "94.46"
$ws.Range("J11").Value = "0.00"
$ws.Range("K11").Value = "0.0"
$ws.Range("L11").Value = "921279283"
$ws.Range("N11").Value = "NO"
$ws.Range("O11").Value = "NO"
$ws.Range("P11").Value = "NO"
$ws.Range("Q11").Value = "NO"

# Row 12
$ws.Range("A12").Value = "''U91200304"
$ws.Range("B12").Value = "''4682412800137"
$ws.Range("C12").Value = "''240516090004339"
$ws.Range("D12").Value = "SALTOS"
$ws.Range("E12").Value = "ELENI"
$ws.Range("F12").Value = "02/16/2024"
$ws.Range("G12").Value = "800.00"
$ws.Range("H12").Value = "65.03"
$ws.Range("I12").Value = "92.90"
$ws.Range("J12").Value = "0.00"
$ws.Range("K12").Value = "0.0"
$ws.Range("L12").Value = "921279283"
$ws.Range("N12").Value = "NO"
$ws.Range("O12").Value = "NO"
$ws.Range("P12").Value = "NO"
$ws.Range("Q12").Value = "NO"

# Row 13
$ws.Range("A13").Value = "''U91200304"
$ws.Range("B13").Value = "''4682412800138"
$ws.Range("C13").Value = "''240516090004339"
$ws.Range("D13").Value = "SALTOS"
$ws.Range("E13").Value = "ELENI"
$ws.Range("F13").Value = "02/21/2024"
$ws.Range("G13").Value = "800.00"
$ws.Range("H13").Value = "65.03"
$ws.Range("I13").Value = "92.90"
$ws.Range("J13").Value = "0.00"
$ws.Range("K13").Value = "0.0"
$ws.Range("L13").Value = "921279283"
$ws.Range("N13").Value = "NO"
$ws.Range("O13").Value = "NO"
$ws.Range("P13").Value = "NO"
$ws.Range("Q13").Value = "NO"

# Row 14
$ws.Range("A14").Value = "''U91200304"
$ws.Range("B14").Value = "''4682412800139"
$ws.Range("C14").Value = "''240516090004339"
$ws.Range("D14").Value = "SALTOS"
$ws.Range("E14").Value = "ELENI"
$ws.Range("F14").Value = "03/20/2024"
$ws.Range("G14").Value = "800.00"
$ws.Range("H14").Value = "66.12"
$ws.Range("I14").Value = "94.46"
$ws.Range("J14").Value = "0.00"
$ws.Range("K14").Value = "0.0"
$ws.Range("L14").Value = "921279283"
$ws.Range("N14").Value = "NO"
$ws.Range("O14").Value = "NO"
$ws.Range("P14").Value = "NO"
$ws.Range("Q14").Value = "NO"

# Row 15
$ws.Range("A15").Value = "''U91200304"
$ws.Range("B15").Value = "''4682412800140"
$ws.Range("C15").Value = "''240516090004339"
$ws.Range("D15").Value = "SALTOS"
$ws.Range("E15").Value = "ELENI"
$ws.Range("F15").Value = "04/10/2024"
$ws.Range("G15").Value = "800.00"
$ws.Range("H15").Value = "66.12"
$ws.Range("I15").Value = "94.46"
$ws.Range("J15").Value = "0.00"
$ws.Range("K15").Value = "0.0"
$ws.Range("L15").Value = "921279283"
$ws.Range("N15").Value = "NO"
$ws.Range("O15").Value = "NO"
$ws.Range("P15").Value = "NO"
$ws.Range("Q15").Value = "NO"

# Remove the two trailing rows that are no longer part of the export
$ws.Rows("16:17").Delete()
